{"js": "// Auto-generated: replace each equation's old text with its new text.\n// Each (old, new) pair is unique and non-overlapping, so direct search+replace is safe.\nconst pairs = [\n  [\"67-49=18\", \"5+77=82\"],\n  [\"38+40=78\", \"19-9=10\"],\n  [\"47+15=62\", \"49+20=69\"],\n  [\"76+4=80\", \"68-45=23\"],\n  [\"58+2=60\", \"39-8=31\"],\n  [\"4+46=50\", \"80-64=16\"],\n  [\"31-21=10\", \"18+58=76\"],\n  [\"50+21=71\", \"94-53=41\"],\n  [\"78-67=11\", \"21+28=49\"],\n  [\"55+28=83\", \"68-23=45\"],\n  [\"11+71=82\", \"28+17=45\"],\n  [\"44-14=30\", \"70-2=68\"],\n  [\"36+56=92\", \"31+28=59\"],\n  [\"4+6=10\", \"13+51=64\"],\n  [\"36+44=80\", \"78+3=81\"],\n  [\"89-62=27\", \"93-72=21\"],\n  [\"72-57=15\", \"63-42=21\"],\n  [\"24+32=56\", \"55-5=50\"],\n  [\"93-19=74\", \"8+51=59\"],\n  [\"56-38=18\", \"82-66=16\"],\n  [\"29+64=93\", \"6+80=86\"],\n  [\"30+10=40\", \"42+19=61\"],\n  [\"37+55=92\", \"96-63=33\"],\n  [\"66+29=95\", \"11+68=79\"],\n  [\"97-50=47\", \"1+30=31\"],\n  [\"92-54=38\", \"98-7=91\"],\n  [\"23+46=69\", \"11+80=91\"],\n  [\"6+92=98\", \"10+0=10\"],\n  [\"47+12=59\", \"40-25=15\"],\n  [\"91-19=72\", \"38-0=38\"],\n  [\"81-62=19\", \"49+21=70\"],\n  [\"10-1=9\", \"1+3=4\"],\n  [\"26+64=90\", \"79-35=44\"],\n  [\"86-4=82\", \"1+10=11\"],\n  [\"91+5=96\", \"34-15=19\"],\n  [\"54-42=12\", \"87-15=72\"],\n  [\"20+15=35\", \"70+0=70\"],\n  [\"6+21=27\", \"82-48=34\"],\n  [\"3+95=98\", \"48-21=27\"],\n  [\"57-55=2\", \"79-8=71\"],\n  [\"67-7=60\", \"54-35=19\"],\n  [\"88+0=88\", \"82-67=15\"],\n  [\"0+56=56\", \"86-32=54\"],\n  [\"53-38=15\", \"86-45=41\"],\n  [\"97-64=33\", \"77-31=46\"],\n  [\"31-28=3\", \"42+56=98\"],\n  [\"53+26=79\", \"99-25=74\"],\n  [\"70-9=61\", \"59-57=2\"],\n  [\"94+4=98\", \"52-46=6\"],\n  [\"29-5=24\", \"56+20=76\"],\n  [\"98-0=98\", \"34+47=81\"],\n  [\"95-40=55\", \"79-32=47\"],\n  [\"66-38=28\", \"49+28=77\"],\n  [\"69-54=15\", \"64+15=79\"],\n  [\"76-19=57\", \"1+52=53\"],\n  [\"61+16=77\", \"73+21=94\"],\n  [\"52-50=2\", \"51+36=87\"],\n  [\"48-43=5\", \"8+3=11\"],\n  [\"45-21=24\", \"29+37=66\"],\n  [\"49-9=40\", \"41+47=88\"],\n  [\"90-1=89\", \"4+22=26\"],\n  [\"45+20=65\", \"1+4=5\"],\n  [\"10+86=96\", \"79-27=52\"],\n  [\"82+10=92\", \"83-62=21\"],\n  [\"86+6=92\", \"85-46=39\"],\n  [\"37+1=38\", \"22+19=41\"],\n  [\"95-88=7\", \"53-13=40\"],\n  [\"71-57=14\", \"46-6=40\"],\n  [\"46+37=83\", \"89-23=66\"],\n  [\"82-4=78\", \"82+5=87\"],\n  [\"97-25=72\", \"58-56=2\"],\n  [\"96-30=66\", \"88-4=84\"],\n  [\"22+7=29\", \"55-42=13\"],\n  [\"54+20=74\", \"82-43=39\"],\n  [\"80-67=13\", \"21+28=49\"],\n  [\"9-4=5\", \"47+39=86\"],\n  [\"44+33=77\", \"29+39=68\"],\n  [\"91-31=60\", \"23+13=36\"],\n  [\"8+68=76\", \"54-44=10\"],\n  [\"17+43=60\", \"29+68=97\"],\n  [\"52-32=20\", \"13-12=1\"],\n  [\"16+30=46\", \"1+78=79\"],\n  [\"42+2=44\", \"98-77=21\"],\n  [\"78-24=54\", \"52-51=1\"],\n  [\"28+59=87\", \"85-76=9\"],\n  [\"48-32=16\", \"53+17=70\"],\n  [\"50-45=5\", \"95-57=38\"],\n  [\"75+18=93\", \"38-22=16\"],\n  [\"80-13=67\", \"74-30=44\"],\n  [\"63+22=85\", \"31+23=54\"],\n  [\"10+68=78\", \"20+29=49\"],\n  [\"22+8=30\", \"34+48=82\"],\n  [\"47-0=47\", \"22+72=94\"],\n  [\"75-68=7\", \"49+18=67\"],\n  [\"69-32=37\", \"31-8=23\"],\n  [\"42-29=13\", \"88-84=4\"],\n  [\"85-1=84\", \"62+20=82\"],\n  [\"53+23=76\", \"3+9=12\"],\n  [\"65+27=92\", \"13+15=28\"],\n  [\"28+23=51\", \"2+4=6\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "# Auto-generated: replace each equation's old text with its new text using Find/Replace.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"67-49=18\", \"5+77=82\")\n    ,@(\"38+40=78\", \"19-9=10\")\n    ,@(\"47+15=62\", \"49+20=69\")\n    ,@(\"76+4=80\", \"68-45=23\")\n    ,@(\"58+2=60\", \"39-8=31\")\n    ,@(\"4+46=50\", \"80-64=16\")\n    ,@(\"31-21=10\", \"18+58=76\")\n    ,@(\"50+21=71\", \"94-53=41\")\n    ,@(\"78-67=11\", \"21+28=49\")\n    ,@(\"55+28=83\", \"68-23=45\")\n    ,@(\"11+71=82\", \"28+17=45\")\n    ,@(\"44-14=30\", \"70-2=68\")\n    ,@(\"36+56=92\", \"31+28=59\")\n    ,@(\"4+6=10\", \"13+51=64\")\n    ,@(\"36+44=80\", \"78+3=81\")\n    ,@(\"89-62=27\", \"93-72=21\")\n    ,@(\"72-57=15\", \"63-42=21\")\n    ,@(\"24+32=56\", \"55-5=50\")\n    ,@(\"93-19=74\", \"8+51=59\")\n    ,@(\"56-38=18\", \"82-66=16\")\n    ,@(\"29+64=93\", \"6+80=86\")\n    ,@(\"30+10=40\", \"42+19=61\")\n    ,@(\"37+55=92\", \"96-63=33\")\n    ,@(\"66+29=95\", \"11+68=79\")\n    ,@(\"97-50=47\", \"1+30=31\")\n    ,@(\"92-54=38\", \"98-7=91\")\n    ,@(\"23+46=69\", \"11+80=91\")\n    ,@(\"6+92=98\", \"10+0=10\")\n    ,@(\"47+12=59\", \"40-25=15\")\n    ,@(\"91-19=72\", \"38-0=38\")\n    ,@(\"81-62=19\", \"49+21=70\")\n    ,@(\"10-1=9\", \"1+3=4\")\n    ,@(\"26+64=90\", \"79-35=44\")\n    ,@(\"86-4=82\", \"1+10=11\")\n    ,@(\"91+5=96\", \"34-15=19\")\n    ,@(\"54-42=12\", \"87-15=72\")\n    ,@(\"20+15=35\", \"70+0=70\")\n    ,@(\"6+21=27\", \"82-48=34\")\n    ,@(\"3+95=98\", \"48-21=27\")\n    ,@(\"57-55=2\", \"79-8=71\")\n    ,@(\"67-7=60\", \"54-35=19\")\n    ,@(\"88+0=88\", \"82-67=15\")\n    ,@(\"0+56=56\", \"86-32=54\")\n    ,@(\"53-38=15\", \"86-45=41\")\n    ,@(\"97-64=33\", \"77-31=46\")\n    ,@(\"31-28=3\", \"42+56=98\")\n    ,@(\"53+26=79\", \"99-25=74\")\n    ,@(\"70-9=61\", \"59-57=2\")\n    ,@(\"94+4=98\", \"52-46=6\")\n    ,@(\"29-5=24\", \"56+20=76\")\n    ,@(\"98-0=98\", \"34+47=81\")\n    ,@(\"95-40=55\", \"79-32=47\")\n    ,@(\"66-38=28\", \"49+28=77\")\n    ,@(\"69-54=15\", \"64+15=79\")\n    ,@(\"76-19=57\", \"1+52=53\")\n    ,@(\"61+16=77\", \"73+21=94\")\n    ,@(\"52-50=2\", \"51+36=87\")\n    ,@(\"48-43=5\", \"8+3=11\")\n    ,@(\"45-21=24\", \"29+37=66\")\n    ,@(\"49-9=40\", \"41+47=88\")\n    ,@(\"90-1=89\", \"4+22=26\")\n    ,@(\"45+20=65\", \"1+4=5\")\n    ,@(\"10+86=96\", \"79-27=52\")\n    ,@(\"82+10=92\", \"83-62=21\")\n    ,@(\"86+6=92\", \"85-46=39\")\n    ,@(\"37+1=38\", \"22+19=41\")\n    ,@(\"95-88=7\", \"53-13=40\")\n    ,@(\"71-57=14\", \"46-6=40\")\n    ,@(\"46+37=83\", \"89-23=66\")\n    ,@(\"82-4=78\", \"82+5=87\")\n    ,@(\"97-25=72\", \"58-56=2\")\n    ,@(\"96-30=66\", \"88-4=84\")\n    ,@(\"22+7=29\", \"55-42=13\")\n    ,@(\"54+20=74\", \"82-43=39\")\n    ,@(\"80-67=13\", \"21+28=49\")\n    ,@(\"9-4=5\", \"47+39=86\")\n    ,@(\"44+33=77\", \"29+39=68\")\n    ,@(\"91-31=60\", \"23+13=36\")\n    ,@(\"8+68=76\", \"54-44=10\")\n    ,@(\"17+43=60\", \"29+68=97\")\n    ,@(\"52-32=20\", \"13-12=1\")\n    ,@(\"16+30=46\", \"1+78=79\")\n    ,@(\"42+2=44\", \"98-77=21\")\n    ,@(\"78-24=54\", \"52-51=1\")\n    ,@(\"28+59=87\", \"85-76=9\")\n    ,@(\"48-32=16\", \"53+17=70\")\n    ,@(\"50-45=5\", \"95-57=38\")\n    ,@(\"75+18=93\", \"38-22=16\")\n    ,@(\"80-13=67\", \"74-30=44\")\n    ,@(\"63+22=85\", \"31+23=54\")\n    ,@(\"10+68=78\", \"20+29=49\")\n    ,@(\"22+8=30\", \"34+48=82\")\n    ,@(\"47-0=47\", \"22+72=94\")\n    ,@(\"75-68=7\", \"49+18=67\")\n    ,@(\"69-32=37\", \"31-8=23\")\n    ,@(\"42-29=13\", \"88-84=4\")\n    ,@(\"85-1=84\", \"62+20=82\")\n    ,@(\"53+23=76\", \"3+9=12\")\n    ,@(\"65+27=92\", \"13+15=28\")\n    ,@(\"28+23=51\", \"2+4=6\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
